$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their literal text formatting (e.g. trailing
# zeros / multi-dot thousand separators) instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.268.24"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "1.860.78"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "0.7035"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").Value = "237.71"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.08219"
$ws.Range("E8").Value = "  +9.44%  "
$ws.Range("D9").Value = "0.3041"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").Value = "23.27"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("D11").Value = "0.08179"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").Value = "1.873.93"
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("D13").Value = "0.7161"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").Value = "5.177"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").Value = "89.20"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").Value = "29.289.62"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "5.777"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("E18").Value = "  +2.42%  "
$ws.Range("D19").Value = "0.000007848"
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("D20").Value = "237.40"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "2.109.39"
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "7.461"
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("D25").Value = "162.04"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").Value = "8.983"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "0.1442"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").Value = "1.968"
$ws.Range("E29").Value = "  +1.64%  "
$ws.Range("D30").Value = "1.440"
$ws.Range("E30").Value = "  +3.99%  "
$ws.Range("D31").Value = "4.429"
$ws.Range("E31").Value = "  -3.01%  "
$ws.Range("D32").Value = "1.482"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("E33").Value = "  +1.50%  "
$ws.Range("D34").Value = "0.05216"
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("D35").Value = "1.171"
$ws.Range("E35").Value = "  -1.26%  "
$ws.Range("D36").Value = "0.7083"
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("D37").Value = "1.002"
$ws.Range("E37").Value = "  -3.76%  "
$ws.Range("D38").Value = "2.668"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("D39").Value = "0.01851"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").Value = "2.721"
$ws.Range("E40").Value = "  +1.68%  "
$ws.Range("D41").Value = "1.139.97"
$ws.Range("E41").Value = "  +5.88%  "
$ws.Range("D42").Value = "0.9172"
$ws.Range("E42").Value = "  -3.32%  "
$ws.Range("D43").Value = "5.966"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").Value = "0.4284"
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("D45").Value = "70.84"
$ws.Range("E45").Value = "  +1.12%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "102.69"
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("D49").Value = "2.007.64"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").Value = "9.183"
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("D51").Value = "6.972"
$ws.Range("E51").Value = "  -1.28%  "
